# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet right before "总计", populated with the
#    per-fund holdings table for that quarter.
# 2. Prepend a "2022-Q1" summary row to the "总计" sheet, shifting the
#    existing quarterly summary rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q1" worksheet positioned right before 总计
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet handles in this host are position-bound, so once Add()
# shifts "总计" over by one slot, the original $totalSheet variable would
# silently alias the brand-new "2022-Q1" sheet instead. Re-resolve "总计"
# by name so Step 2 below edits the correct sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Match the other quarter sheets' outline + page-margin setup (same
# "Normal" Excel defaults used throughout this workbook: 0.75"/1"/0.5",
# PageSetup's margin properties are expressed in points, hence *72).
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newPs = $newSheet.PageSetup
$newPs.LeftMargin = 0.75 * 72
$newPs.RightMargin = 0.75 * 72
$newPs.TopMargin = 1 * 72
$newPs.BottomMargin = 1 * 72
$newPs.HeaderMargin = 0.5 * 72
$newPs.FooterMargin = 0.5 * 72

# Copy the header row formatting/labels (columns B:H) from an existing
# quarter sheet - they all share the same "基金代码 / 基金名称 / ..." header.
$template.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Copy the bold/bordered index-column style used for column A (A2) down
# across every data row so the whole A2:A23 range is pre-formatted.
$template.Range("A2").Copy($newSheet.Range("A2:A23"))

# Per-fund holdings data for 2022-Q1: code, name, scale, total stock
# position, position ratio, held market value (亿元), position rank.
$fundData = @(
    @('519196','万家新兴蓝筹灵活配置混合','21.26','80.70','4.73','1.0056',9),
    @('011177','博时汇融回报一年持有期混合A','39.26','65.35','2.16','0.8480',10),
    @('519195','万家品质生活灵活配置混合','17.66','79.87','4.32','0.7629',9),
    @('005094','万家臻选混合','13.43','73.12','5.23','0.7024',7),
    @('519181','万家和谐增长混合','11.51','88.62','4.84','0.5571',8),
    @('161912','万家社会责任18个月定期开放混合（LOF）A','13.56','88.11','4.06','0.5505',10),
    @('006864','国联安核心资产策略混合','5.64','91.69','5.41','0.3051',4),
    @('121006','国投瑞银稳健增长混合','7.34','65.10','2.20','0.1615',5),
    @('161232','国投瑞银瑞盛灵活配置混合','4.28','94.52','3.66','0.1566',10),
    @('000663','国投瑞银美丽中国灵活配置混合','3.55','92.92','3.91','0.1388',10),
    @('161225','国投瑞银瑞盈灵活配置混合（LOF）','2.28','94.46','4.03','0.0919',5),
    @('519198','万家颐和灵活配置混合','1.78','91.03','4.69','0.0835',10),
    @('161233','国投瑞银瑞泰多策略灵活配置混合（LOF）','7.05','24.66','0.94','0.0663',8),
    @('005117','金信价值精选灵活配置混合A','0.83','86.42','5.13','0.0426',4),
    @('002885','摩根士丹利华鑫万众创新灵活配置混合','0.86','94.36','4.70','0.0404',9),
    @('009317','金信核心竞争力灵活配置混合','0.19','89.48','9.80','0.0186',3),
    @('161913','万家社会责任18个月定期开放混合（LOF）C','0.44','88.11','4.06','0.0179',10),
    @('011178','博时汇融回报一年持有期混合C','0.16','65.35','2.16','0.0035',10),
    @('162107','金鹰量化精选股票（LOF）','0.06','93.77','5.64','0.0034',8),
    @('005118','金信价值精选灵活配置混合C','0.05','86.42','5.13','0.0026',4),
    @('010605','创金合信鑫祥混合A','0.03','29.57','1.10','0.0003',4),
    @('010606','创金合信鑫祥混合C','0.02','29.57','1.10','0.0002',4)
)

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = $i + 2
    $rec = $fundData[$i]

    $newSheet.Cells.Item($row, 1).Value = $i

    # Force text storage (leading zeros / fixed-decimal strings) with a
    # leading apostrophe, same as typing the value into Excel manually.
    $newSheet.Cells.Item($row, 2).Value = "'" + $rec[0]
    $newSheet.Cells.Item($row, 3).Value = $rec[1]
    $newSheet.Cells.Item($row, 4).Value = "'" + $rec[2]
    $newSheet.Cells.Item($row, 5).Value = "'" + $rec[3]
    $newSheet.Cells.Item($row, 6).Value = "'" + $rec[4]
    $newSheet.Cells.Item($row, 7).Value = "'" + $rec[5]
    $newSheet.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row to 总计, shifting older rows down
# ---------------------------------------------------------------------

# Shift existing summary rows 2-6 down to 3-7, preserving their formatting.
$totalSheet.Range("A2:D6").Copy($totalSheet.Range("A3:D7"))

# Renumber the shifted rows' index column (A) to stay sequential.
for ($r = 3; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Write the new 2022-Q1 summary into the now-vacated row 2.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 22
$totalSheet.Cells.Item(2, 4).Value = 5.56
